$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.898.17"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").Value = "1.625.38"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.44"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +8.19%  "

$ws.Range("E9").Value = "  +2.63%  "

$ws.Range("E10").Value = "  +1.76%  "

$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("D12").Value = "1.858.73"
$ws.Range("E12").Value = "  +1.23%  "

$ws.Range("D13").Value = "1.627.19"
$ws.Range("E13").Value = "  +1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.572"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.68%  "

$ws.Range("E15").Value = "  +4.63%  "

$ws.Range("D16").Value = "29.954.17"
$ws.Range("E16").Value = "  +0.89%  "

$ws.Range("E17").Value = "  +17.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("D20").Value = "0.0₃0705"
$ws.Range("E20").Value = "  +1.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.62"
$ws.Range("D23").ClearFormats()

$ws.Range("E24").Value = "  +1.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.65"
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.68"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.12%  "

$ws.Range("E27").Value = "  +2.35%  "

$ws.Range("E28").Value = "  +2.96%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  +2.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.22"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.84%  "

$ws.Range("D34").Value = "1.422.02"
$ws.Range("E34").Value = "  -0.68%  "

$ws.Range("E35").Value = "  +6.48%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.30"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("E39").Value = "  +3.35%  "

$ws.Range("E41").Value = "  +0.87%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.831"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.75%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0500"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "54.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.33"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.30%  "

$ws.Range("E46").Value = "  +8.80%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").Value = "  +1.97%  "

$ws.Range("D49").Value = "1.766.32"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.38%  "

$ws.Range("E51").Value = "  +6.49%  "
